# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
# Both sheets mirror the same event list, so the same row -> new-value map
# applies to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 62
    5  = 62
    6  = 12536
    7  = 12536
    10 = 503
    12 = 1143
    13 = 938
    14 = 13643
    15 = 13944
    17 = 166
    20 = 1048
    23 = 467
    24 = 5039
    25 = 248
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
